$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these ID cells remain text (avoid numeric auto-conversion)
$ws.Range("B234").NumberFormat = "@"
$ws.Range("B235").NumberFormat = "@"
$ws.Range("B236").NumberFormat = "@"
$ws.Range("B237").NumberFormat = "@"
$ws.Range("B238").NumberFormat = "@"
$ws.Range("B239").NumberFormat = "@"
$ws.Range("B240").NumberFormat = "@"

# Row 114
$ws.Range("B114").Value = 7559469
$ws.Range("E114").Value = "Montevideo Wanderers"
$ws.Range("F114").Value = "Penarol"
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = "D"
$ws.Range("L114").Value = 4.75
$ws.Range("M114").Value = 3.4
$ws.Range("N114").Value = 1.7
$ws.Range("O114").Value = 2.7
$ws.Range("Q114").Value = 2.45
$ws.Range("R114").Value = 0
$ws.Range("S114").Value = 2.05
$ws.Range("T114").Value = 1.8
$ws.Range("U114").Value = 2.5
$ws.Range("V114").Value = 1.975
$ws.Range("W114").Value = 1.875
$ws.Range("X114").Value = -1
$ws.Range("Y114").Value = 2.2
$ws.Range("AA114").Value = 0
$ws.Range("AB114").Value = 0
$ws.Range("AC114").Value = -1
$ws.Range("AD114").Value = 0.875

# Row 115
$ws.Range("B115").Value = 7559468
$ws.Range("E115").Value = "Liverpool Montevideo"
$ws.Range("F115").Value = "CA River Plate"
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 1
$ws.Range("I115").Value = 1
$ws.Range("K115").Value = "H"
$ws.Range("L115").Value = 1.7
$ws.Range("M115").Value = 3
$ws.Range("N115").Value = 5.75
$ws.Range("O115").Value = 1.833
$ws.Range("Q115").Value = 4.5
$ws.Range("R115").Value = -0.5
$ws.Range("S115").Value = 1.925
$ws.Range("T115").Value = 1.925
$ws.Range("U115").Value = 2.25
$ws.Range("V115").Value = 2.025
$ws.Range("W115").Value = 1.825
$ws.Range("X115").Value = 0.833
$ws.Range("Y115").Value = -1
$ws.Range("AA115").Value = 0.925
$ws.Range("AB115").Value = -1
$ws.Range("AC115").Value = 1.025
$ws.Range("AD115").Value = -1

# Row 117
$ws.Range("B117").Value = 7013886
$ws.Range("E117").Value = "Racing Club de Montevideo"
$ws.Range("F117").Value = "Cerro"
$ws.Range("G117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = "A"
$ws.Range("L117").Value = 2.25
$ws.Range("M117").Value = 3.1
$ws.Range("N117").Value = 3.25
$ws.Range("O117").Value = 2.25
$ws.Range("P117").Value = 2.875
$ws.Range("Q117").Value = 3.5
$ws.Range("R117").Value = -0.25
$ws.Range("S117").Value = 1.95
$ws.Range("T117").Value = 1.9
$ws.Range("U117").Value = 2
$ws.Range("V117").Value = 1.925
$ws.Range("W117").Value = 1.925
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = 2.5
$ws.Range("AB117").Value = 0.8999999999999999
$ws.Range("AD117").Value = 0.925

# Row 118
$ws.Range("B118").Value = 7013885
$ws.Range("E118").Value = "La Luz"
$ws.Range("F118").Value = "Atletico Fenix Montevideo"
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 3
$ws.Range("M118").Value = 3
$ws.Range("N118").Value = 2.4
$ws.Range("O118").Value = 2.9
$ws.Range("P118").Value = 2.75
$ws.Range("Q118").Value = 2.6
$ws.Range("R118").Value = 0
$ws.Range("S118").Value = 2.025
$ws.Range("T118").Value = 1.825
$ws.Range("U118").Value = 2
$ws.Range("V118").Value = 2.025
$ws.Range("W118").Value = 1.825
$ws.Range("Z118").Value = 1.6
$ws.Range("AB118").Value = 0.825
$ws.Range("AC118").Value = 0
$ws.Range("AD118").Value = 0

# Row 119
$ws.Range("B119").Value = 7013702
$ws.Range("E119").Value = "Defensor Sporting"
$ws.Range("F119").Value = "Danubio"
$ws.Range("J119").Value = 1
$ws.Range("L119").Value = 1.8
$ws.Range("M119").Value = 3.6
$ws.Range("N119").Value = 4.2
$ws.Range("O119").Value = 1.8
$ws.Range("P119").Value = 3.6
$ws.Range("Q119").Value = 4.2
$ws.Range("R119").Value = -0.75
$ws.Range("S119").Value = 2.05
$ws.Range("T119").Value = 1.8
$ws.Range("U119").Value = 2.25
$ws.Range("V119").Value = 1.85
$ws.Range("W119").Value = 2
$ws.Range("Z119").Value = 3.2
$ws.Range("AB119").Value = 0.8
$ws.Range("AC119").Value = -0.5
$ws.Range("AD119").Value = 0.5

# Row 120
$ws.Range("B120").Value = 7013409
$ws.Range("E120").Value = "Nacional De Football"
$ws.Range("F120").Value = "Torque"
$ws.Range("G120").Value = 1
$ws.Range("J120").Value = 1
$ws.Range("K120").Value = "D"
$ws.Range("L120").Value = 1.666
$ws.Range("M120").Value = 3.9
$ws.Range("N120").Value = 4.5
$ws.Range("O120").Value = 1.615
$ws.Range("P120").Value = 4
$ws.Range("Q120").Value = 4.75
$ws.Range("R120").Value = -0.75
$ws.Range("S120").Value = 1.8
$ws.Range("T120").Value = 2.05
$ws.Range("U120").Value = 2.75
$ws.Range("V120").Value = 1.95
$ws.Range("W120").Value = 1.9
$ws.Range("Y120").Value = 3
$ws.Range("Z120").Value = -1
$ws.Range("AB120").Value = 1.05
$ws.Range("AD120").Value = 0.8999999999999999

# Row 234
$ws.Range("B234").Value = "8261936"
$ws.Range("D234").Value = 45444.72916666666
$ws.Range("E234").Value = "Montevideo Wanderers"
$ws.Range("F234").Value = "Rampla Juniors"
$ws.Range("L234").Value = 2.15
$ws.Range("N234").Value = 3.4
$ws.Range("O234").Value = 2.1
$ws.Range("P234").Value = 3.1
$ws.Range("Q234").Value = 3.7
$ws.Range("S234").Value = 1.825
$ws.Range("T234").Value = 2.025
$ws.Range("V234").Value = 2.025
$ws.Range("W234").Value = 1.825

# Row 235
$ws.Range("B235").Value = "8260831"
$ws.Range("D235").Value = 45445.41666666666
$ws.Range("E235").Value = "CA River Plate"
$ws.Range("F235").Value = "Racing Club de Montevideo"
$ws.Range("L235").Value = 2.75
$ws.Range("M235").Value = 3.1
$ws.Range("N235").Value = 2.625
$ws.Range("O235").Value = 3.2
$ws.Range("P235").Value = 3.1
$ws.Range("Q235").Value = 2.3
$ws.Range("R235").Value = 0.25
$ws.Range("S235").Value = 1.85
$ws.Range("T235").Value = 2
$ws.Range("V235").Value = 1.875
$ws.Range("W235").Value = 1.975

# Row 236
$ws.Range("B236").Value = "8260830"
$ws.Range("D236").Value = 45445.625
$ws.Range("E236").Value = "Atletico Fenix Montevideo"
$ws.Range("F236").Value = "Penarol"
$ws.Range("L236").Value = 5.5
$ws.Range("M236").Value = 3.6
$ws.Range("N236").Value = 1.615
$ws.Range("O236").Value = 5.5
$ws.Range("P236").Value = 3.75
$ws.Range("Q236").Value = 1.615
$ws.Range("R236").Value = 0.75
$ws.Range("S236").Value = 2.025
$ws.Range("T236").Value = 1.825
$ws.Range("V236").Value = 1.85
$ws.Range("W236").Value = 2

# Row 237
$ws.Range("B237").Value = "8261937"
$ws.Range("D237").Value = 45445.72916666666
$ws.Range("E237").Value = "Defensor Sporting"
$ws.Range("F237").Value = "Miramar Misiones"
$ws.Range("L237").Value = 1.615
$ws.Range("M237").Value = 3.75
$ws.Range("N237").Value = 5.5
$ws.Range("O237").Value = 1.65
$ws.Range("P237").Value = 3.7
$ws.Range("Q237").Value = 5.25
$ws.Range("R237").Value = -0.75
$ws.Range("S237").Value = 1.85
$ws.Range("T237").Value = 2
$ws.Range("U237").Value = 2.5

# Row 238
$ws.Range("B238").Value = "8260832"
$ws.Range("D238").Value = 45446.52083333334
$ws.Range("E238").Value = "Liverpool Montevideo"
$ws.Range("F238").Value = "Cerro"
$ws.Range("L238").Value = 1.833
$ws.Range("M238").Value = 3.5
$ws.Range("N238").Value = 4.333
$ws.Range("O238").Value = 1.85
$ws.Range("P238").Value = 3.5
$ws.Range("Q238").Value = 4.2
$ws.Range("R238").Value = -0.5
$ws.Range("S238").Value = 1.875
$ws.Range("T238").Value = 1.975
$ws.Range("V238").Value = 2
$ws.Range("W238").Value = 1.85

# Row 239
$ws.Range("B239").Value = "8260833"
$ws.Range("D239").Value = 45446.625
$ws.Range("E239").Value = "Danubio"
$ws.Range("F239").Value = "Deportivo Maldonado"
$ws.Range("L239").Value = 2.375
$ws.Range("M239").Value = 3.2
$ws.Range("N239").Value = 3
$ws.Range("O239").Value = 2.35
$ws.Range("P239").Value = 3.2
$ws.Range("Q239").Value = 3
$ws.Range("R239").Value = -0.25
$ws.Range("S239").Value = 2.05
$ws.Range("T239").Value = 1.8
$ws.Range("U239").Value = 2.25

# Row 240
$ws.Range("B240").Value = "8260829"
$ws.Range("D240").Value = 45447.625
$ws.Range("E240").Value = "Club Atletico Progreso"
$ws.Range("F240").Value = "Cerro Largo"
$ws.Range("L240").Value = 2.25
$ws.Range("N240").Value = 3.2
$ws.Range("O240").Value = 2.3
$ws.Range("Q240").Value = 3.1
$ws.Range("S240").Value = 2.025
$ws.Range("T240").Value = 1.825
$ws.Range("V240").Value = 1.925
$ws.Range("W240").Value = 1.925
